$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tb = $s.Shapes.AddTextbox(1, 0, 0, 100, 50)
$tr = $tb.TextFrame.TextRange
Write-Host "Before: [$($tr.Text)]"
$tr.InsertSlideNumber()
Write-Host "After: [$($tr.Text)]"
Write-Host "Shape count: $($s.Shapes.Count)"
